# Generate Report for Handback
# Adds "Latest Target File" / "Latest Handback File" hyperlinked values to the
# zh-cn and de-de localization-status sheets, flips the status text to the
# "handed back" state, and stamps the handback datetime.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl_5077 = "https://github.com/OpenLocalizationTest/oltest/blob/3e24e39ffa24d9e02c342cff803a03ec918ee08d/e2e/5077eaab-51ab-4868-9300-0c7db760429e.md"
$mdUrl_863e = "https://github.com/OpenLocalizationTest/oltest/blob/3e24e39ffa24d9e02c342cff803a03ec918ee08d/e2e/863e2609-e7a8-4211-a0aa-5603a3d9c989.md"

# Per-locale worksheet data: the handback timestamp for that locale's run and
# the hyperlink targets for the translated (.xlf) files that were handed off.
$locales = @(
    @{
        Sheet = "zh-cn"
        Stamp = "2016-03-11 14:42:53"
        Xlf = @{
            "2" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf07b32ce4295dd5fcc7c858a9d2714aef24aafa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5077eaab-51ab-4868-9300-0c7db760429e.d36459d7f03bf69771da7f073578b52c6259f240.zh-cn.xlf"
            "3" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf07b32ce4295dd5fcc7c858a9d2714aef24aafa/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/863e2609-e7a8-4211-a0aa-5603a3d9c989.b19cb24c14e837881260e26bd3365eddca6c3b0d.zh-cn.xlf"
        }
    },
    @{
        Sheet = "de-de"
        Stamp = "2016-03-11 14:42:59"
        Xlf = @{
            "2" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8be2f45578d0ec0b1d98972fc409be00baa6c282/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5077eaab-51ab-4868-9300-0c7db760429e.d36459d7f03bf69771da7f073578b52c6259f240.de-de.xlf"
            "3" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8be2f45578d0ec0b1d98972fc409be00baa6c282/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/863e2609-e7a8-4211-a0aa-5603a3d9c989.b19cb24c14e837881260e26bd3365eddca6c3b0d.de-de.xlf"
        }
    }
)

$mdUrlByRow = @{ "2" = $mdUrl_5077; "3" = $mdUrl_863e }

# The Overview sheet mirrors the same per-language status text (it showed
# "Ready for handoff" for every row/language before this handback run).
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusText
$ov.Range("C2").Value = $statusText
$ov.Range("B3").Value = $statusText
$ov.Range("C3").Value = $statusText

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)
    $stamp = $loc.Stamp

    for ($row = 2; $row -le 3; $row++) {
        $aCell = $ws.Range("A$row")
        $dCell = $ws.Range("D$row")
        $fCell = $ws.Range("F$row")
        $gCell = $ws.Range("G$row")

        $aText = $aCell.Value2
        $dText = $dCell.Value2

        $aTarget = $mdUrlByRow["$row"]
        $dTarget = $loc.Xlf["$row"]

        # Latest Target File == same file that was last handed off.
        $ws.Hyperlinks.Add($fCell, $aTarget, [Type]::Missing, [Type]::Missing, $aText) | Out-Null
        $fCell.Style = "Hyperlink"

        # Latest Handback File == the translated file just handed back.
        $ws.Hyperlinks.Add($gCell, $dTarget, [Type]::Missing, [Type]::Missing, $dText) | Out-Null
        $gCell.Style = "Hyperlink"

        # Status flips from "Ready for handoff" to the synced/handed-back state.
        $ws.Range("C$row").Value = $statusText

        # Latest Handback DateTime stamp for this locale's run.
        $ws.Range("H$row").Value = $stamp
    }
}
